# Update the "想去人数" (F column) figures on the 展览, 演出 and 全部类型 sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 446
$ws1.Range("F3").Value  = 1834
$ws1.Range("F5").Value  = 137
$ws1.Range("F6").Value  = 1722
$ws1.Range("F8").Value  = 140
$ws1.Range("F9").Value  = 650
$ws1.Range("F10").Value = 28
$ws1.Range("F11").Value = 59
$ws1.Range("F12").Value = 549
$ws1.Range("F20").Value = 4505
$ws1.Range("F21").Value = 37
$ws1.Range("F22").Value = 805
$ws1.Range("F24").Value = 2152
$ws1.Range("F26").Value = 3
$ws1.Range("F27").Value = 2021

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 70

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 446
$ws4.Range("F3").Value  = 1834
$ws4.Range("F5").Value  = 137
$ws4.Range("F6").Value  = 1722
$ws4.Range("F8").Value  = 140
$ws4.Range("F9").Value  = 650
$ws4.Range("F10").Value = 28
$ws4.Range("F11").Value = 59
$ws4.Range("F12").Value = 549
$ws4.Range("F20").Value = 4505
$ws4.Range("F21").Value = 70
$ws4.Range("F22").Value = 37
$ws4.Range("F24").Value = 805
$ws4.Range("F26").Value = 2152
$ws4.Range("F28").Value = 3
$ws4.Range("F29").Value = 2021
